$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds the "Price" values and column E the "Volume(1h)" values.
# All of these cells are stored as plain text in the workbook (t="inlineStr"),
# so we force a text number format before assigning, then restore the default
# "Normal" style so no stray style index gets attached to the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.401.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.060.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.42%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -1.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.059.61"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.451"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000237"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.566.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.358.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.064.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "485.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.708"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("E24").Value = "  +4.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  +1.96%  "
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0823"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.88%  "
$ws.Range("E36").Value = "  +1.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("E39").Value = "  +1.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "440.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("E43").Value = "  +2.41%  "
$ws.Range("E44").Value = "  +4.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0364"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.841.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.70%  "
$ws.Range("E51").Value = "  -0.54%  "
